$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 211, shifting all existing
# data (old rows 211-319) down to rows 213-321. Excel copies the row-211
# formatting (incl. the date NumberFormat already applied to column D)
# down onto the two freshly inserted rows, same as it does for row 213.
$ws.Rows("211:212").Insert()

# New row 211: Chino / Primera, 2022-09-02 ($/caja 10 kilos)
$ws.Range("A211").Value = 8
$ws.Range("B211").Value = "Terminal La Palmera de La Serena"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44806
$ws.Range("E211").Value = 4
$ws.Range("F211").Value = 100112003
$ws.Range("G211").Value = "Ajo"
$ws.Range("H211").Value = "Chino"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 600
$ws.Range("K211").Value = 23500
$ws.Range("L211").Value = 24000
$ws.Range("M211").Value = 23750
$ws.Range("N211").Value = "$/caja 10 kilos"
$ws.Range("O211").Value = "China"
$ws.Range("P211").Value = 2375
$ws.Range("Q211").Value = 10
$ws.Range("R211").Value = "Hortaliza"

# New row 212: Chino / Primera, 2022-09-02 ($/malla 10 kilos)
$ws.Range("A212").Value = 8
$ws.Range("B212").Value = "Terminal La Palmera de La Serena"
$ws.Range("C212").Value = "Coquimbo"
$ws.Range("D212").Value = 44806
$ws.Range("E212").Value = 4
$ws.Range("F212").Value = 100112003
$ws.Range("G212").Value = "Ajo"
$ws.Range("H212").Value = "Chino"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 500
$ws.Range("K212").Value = 25500
$ws.Range("L212").Value = 26000
$ws.Range("M212").Value = 25750
$ws.Range("N212").Value = "$/malla 10 kilos"
$ws.Range("O212").Value = "China"
$ws.Range("P212").Value = 2575
$ws.Range("Q212").Value = 10
$ws.Range("R212").Value = "Hortaliza"
